$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NEW")

# The automated map update removed two reclamation records from the table
# (dimension shrinks from A1:R75 to A1:R73), shifting every following row
# up:
#   - Caso -404 ("Amenabar 3048", row 22)
#   - Caso -575 ("Amenabar 3064"), originally row 56 -- after the first
#     row is deleted everything shifts up by one, so it is now row 55.
$ws.Range("A22").EntireRow.Delete()
$ws.Range("A55").EntireRow.Delete()
